# Split the single run of text in the Title, Author and Abstract
# paragraphs into one run per word plus one run per inter-word space,
# leaving the visible text (and every other paragraph) unchanged.
#
# Word's object model does not expose a direct "split this run here"
# call, and simply assigning Range.Text / toggling a character-format
# property either performs a no-op merge or leaves stray formatting
# markup behind. Temporarily inserting a paragraph break at each word
# boundary *does* cleanly terminate a run, and deleting that same
# paragraph mark again merges the two paragraphs back into one while
# keeping the runs that were on either side of the mark separate - so
# that round trip is used here purely as a "run splitter".

function Split-ParagraphIntoWordRuns {
    param($Document, $ParagraphIndex)

    $para = $Document.Paragraphs($ParagraphIndex)
    $pStart = $para.Range.Start
    $pEnd = $para.Range.End

    # Text of the paragraph, excluding its trailing paragraph mark.
    $fullText = $Document.Range($pStart, $pEnd - 1).Text

    # Tokenize into runs of non-whitespace and runs of whitespace,
    # preserving order - e.g. "Questions: Further" ->
    # "Questions:", " ", "Further".
    $tokens = [System.Text.RegularExpressions.Regex]::Matches($fullText, '\S+|\s+')

    if ($tokens.Count -le 1) {
        return
    }

    # Absolute document offsets of each internal word/space boundary,
    # in ascending order.
    $boundaries = New-Object System.Collections.Generic.List[int]
    $acc = $pStart
    for ($i = 0; $i -lt $tokens.Count - 1; $i++) {
        $acc = $acc + $tokens[$i].Length
        $boundaries.Add($acc)
    }

    # Insert a paragraph break at each boundary. Processing from the
    # highest offset down to the lowest means every insertion happens
    # strictly after all the boundaries still to be processed, so
    # their (not yet shifted) original offsets stay valid to insert
    # at.
    for ($i = $boundaries.Count - 1; $i -ge 0; $i--) {
        $pos = $boundaries[$i]
        $Document.Range($pos, $pos).InsertParagraphAfter()
    }

    # Each of those inserted paragraph marks shifted every boundary
    # after it by one character, so the mark that was inserted for
    # the i-th boundary (0-based, ascending) now actually sits at
    # boundary[i] + i. Delete them starting with the last one so the
    # still-to-be-deleted marks keep their computed offsets valid;
    # this merges all the paragraphs back into the original single
    # paragraph while leaving the runs split apart.
    for ($i = $boundaries.Count - 1; $i -ge 0; $i--) {
        $pos = $boundaries[$i] + $i
        $Document.Range($pos, $pos + 1).Delete()
    }
}

$d = $word.ActiveDocument

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $styleName = $d.Paragraphs($i).Style.NameLocal
    if ($styleName -eq "Title" -or $styleName -eq "Author" -or $styleName -eq "Abstract") {
        Split-ParagraphIntoWordRuns $d $i
    }
}
